$d = $word.ActiveDocument

function Set-ParaXml {
    param($ParaIndex, $Fragment)
    $p = $d.Paragraphs.Item($ParaIndex)
    $r = $p.Range
    $pkg = "<?xml version=`"1.0`" standalone=`"yes`"?><?mso-application progid=`"Word.Document`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:body><w:p>$Fragment</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
    $null = $r.InsertXML($pkg)
}

$frag2 = "<w:r><w:t>*</w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>Genome</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> architecture &lt;-&gt; </w:t></w:r><w:r w:rsidRPr=`"005F53AF`"><w:rPr><w:lang w:val=`"en-US`"/></w:rPr><w:t>gene</w:t></w:r><w:r><w:t xml:space=`"preserve`"> expression </w:t></w:r>"
Set-ParaXml 2 $frag2

$frag3 = "<w:r><w:sym w:font=`"Wingdings`" w:char=`"F0E0`"/></w:r><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:proofErr w:type=`"gramStart`"/><w:r><w:t>some</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:proofErr w:type=`"gramEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>lincRNAs</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>regulate</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> expression cis-</w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>trans</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`">. </w:t></w:r>"
Set-ParaXml 3 $frag3

$frag4 = "<w:r><w:t>*</w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>Explain</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>TADs</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`">, TAD </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>boundaries</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`">, (CTCF </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>enrichment</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`">). </w:t></w:r>"
Set-ParaXml 4 $frag4

$frag5 = "<w:r><w:t xml:space=`"preserve`">*Blabla </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>investigate</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> TAD-</w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>bound</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>lincRNAs</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r w:rsidR=`"005F53AF`"><w:t>properties</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>"
Set-ParaXml 5 $frag5

$frag6 = "<w:proofErr w:type=`"gramStart`"/><w:r><w:t>?TR</w:t></w:r><w:proofErr w:type=`"gramEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>lincRNAs</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t>, impact on traits</w:t></w:r>"
Set-ParaXml 6 $frag6

$frag8 = "<w:proofErr w:type=`"spellStart`"/><w:r><w:t>Results</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>"
Set-ParaXml 8 $frag8

$frag10 = "<w:r><w:t xml:space=`"preserve`">TAD </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>boundaries</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>definition</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>"
Set-ParaXml 10 $frag10

$frag11 = "<w:proofErr w:type=`"spellStart`"/><w:r><w:t>Diff</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>between</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>TADbound</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> and non-</w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>TADbound</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>lincRNAs</w:t></w:r><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/><w:proofErr w:type=`"spellEnd`"/>"
Set-ParaXml 11 $frag11

$frag13 = "<w:proofErr w:type=`"spellStart`"/><w:r><w:t>Genomic</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> DNA </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>is</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>folded</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> onto </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>itself</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`">, </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>forming</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> compact structures </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>that</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> affect </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>gene</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> expression.  On a large </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>scale</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`">, </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>regions</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>presenting</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> a </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>high</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>degree</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> of compaction are </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>classified</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> as </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>heterochromatin</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>while</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>uncondensed</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>regions</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> are </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>classified</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> as </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>euchromatin</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`">. </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>These</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> are </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>respectively</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>associated</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>with</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>lower</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> and </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>higher</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> expression </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>levels</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`">. On a </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>smaller</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>scale</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`">, areas </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>where</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> DNA-DNA interactions are </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>especially</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>frequent</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> are </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>called</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>topologically</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>associated</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>domains</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> (</w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>TADs</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`">). </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>Those</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>domains</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> are </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>conserve</w:t></w:r><w:r><w:t>d</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>across</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>cell</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>lines</w:t></w:r><w:proofErr w:type=`"spellEnd`"/><w:r><w:t xml:space=`"preserve`"> and </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>contain</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>"
Set-ParaXml 13 $frag13
